$d = $word.ActiveDocument

function Find-ParaIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "$needle*") {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# Change 1: insert a new "Wage Increase" paragraph right before the
# "As a director of a department ..." paragraph under "General - Stories".
# ------------------------------------------------------------------
$targetIdx = Find-ParaIndex("As a director of a department")
$target = $d.Paragraphs.Item($targetIdx)
$target.Range.InsertParagraphBefore()
$wageP = $d.Paragraphs.Item($targetIdx)
$wageP.Style = -1
$wageP.Range.Text = "Wage Increase"
$wageP.Range.Font.Name = "Arial"
$wageP.Range.Font.Size = 10
$wageP.Range.Font.SizeBi = 10

# ------------------------------------------------------------------
# Change 2: the "probation" paragraph -
#   "... cannot confirm the probation with permission ... resides in "
# becomes
#   "... cannot confirm the probation without permission ... resides in."
# split across four runs, and the trailing bookmark (_GoBack) is moved
# into a brand-new empty paragraph that follows it.
# ------------------------------------------------------------------
$probIdx = Find-ParaIndex("As member of the Human Resources, I should be able to create and modify, probation")
$probationPara = $d.Paragraphs.Item($probIdx)
$pStart = $probationPara.Range.Start
$full = $probationPara.Range.Text

# --- insert "out" so "with" -> "without" ---
$idxWith = $full.IndexOf("probation with") + "probation with".Length
$posWith = $pStart + $idxWith
$insertRng = $d.Range($posWith, $posWith)
$insertRng.InsertAfter("out")
$outRng = $d.Range($posWith, $posWith + 3)
$outRng.Font.Name = "Arial"
$outRng.Font.Size = 10
$outRng.Font.SizeBi = 10

# --- change trailing "resides in " -> "resides in." as its own run ---
$probationPara2 = $d.Paragraphs.Item($probIdx)
$pStart2 = $probationPara2.Range.Start
$full2 = $probationPara2.Range.Text
$idxResides = $full2.IndexOf("resides in")
$posResides = $pStart2 + $idxResides
$pEndExclMark = $probationPara2.Range.End - 1
$tailRng = $d.Range($posResides, $pEndExclMark)
$tailRng.Text = "resides in."
$tailRng2 = $d.Range($posResides, $posResides + ("resides in.".Length))
$tailRng2.Font.Name = "Arial"
$tailRng2.Font.Size = 10
$tailRng2.Font.SizeBi = 10

# --- move the _GoBack bookmark out of this paragraph and into a new,
#     following empty paragraph formatted with ind left=720 ---
$probationPara3 = $d.Paragraphs.Item($probIdx)
$probationPara3.Range.InsertParagraphAfter()
$bmParaIdx = $probIdx + 1
$bmPara = $d.Paragraphs.Item($bmParaIdx)

# give the new paragraph some placeholder text so the bookmark range is
# unambiguously anchored inside it, then strip the placeholder back out
$bmPara.Range.InsertAfter("X")
$bmPara2 = $d.Paragraphs.Item($bmParaIdx)
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$rngForBm = $d.Range($bmPara2.Range.Start, $bmPara2.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $rngForBm)
$bmPara3 = $d.Paragraphs.Item($bmParaIdx)
$delRng = $d.Range($bmPara3.Range.Start, $bmPara3.Range.Start + 1)
$delRng.Text = ""

Write-Host "Edit complete"
